$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.310.76"
$ws.Range("E2").Value = "  +1.53%  "

# Row 3
$ws.Range("D3").Value = "3.804.54"
$ws.Range("E3").Value = "  +0.78%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "674.19"
$ws.Range("E5").Value = "  +7.80%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.86"
$ws.Range("E6").Value = "  +1.76%  "

# Row 7
$ws.Range("D7").Value = "3.803.13"
$ws.Range("E7").Value = "  +0.83%  "

# Row 8
$ws.Range("E8").Value = "  +0.10%  "

# Row 9
$ws.Range("E9").Value = "  +0.69%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  +0.78%  "

# Row 11
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.09"
$ws.Range("E11").Value = "  +5.17%  "

# Row 12
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("E12").Value = "  +0.21%  "

# Row 13
$ws.Range("E13").Value = "  -1.63%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.67"
$ws.Range("E14").Value = "  -0.45%  "

# Row 15
$ws.Range("D15").Value = "4.443.22"
$ws.Range("E15").Value = "  +0.65%  "

# Row 16
$ws.Range("D16").Value = "3.802.44"
$ws.Range("E16").Value = "  +0.71%  "

# Row 17
$ws.Range("D17").Value = "70.323.54"
$ws.Range("E17").Value = "  +1.55%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.65"
$ws.Range("E18").Value = "  -0.11%  "

# Row 19
$ws.Range("E19").Value = "  +1.54%  "

# Row 20
$ws.Range("E20").Value = "  +0.49%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.38"
$ws.Range("E21").Value = "  +18.77%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "475.86"
$ws.Range("E22").Value = "  +1.60%  "

# Row 23
$ws.Range("E23").Value = "  +0.88%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.48"
$ws.Range("E24").Value = "  +0.39%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000142"
$ws.Range("E25").Value = "  -4.19%  "

# Row 26
$ws.Range("E26").Value = "  +0.64%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.28"
$ws.Range("E27").Value = "  +2.20%  "

# Row 28
$ws.Range("E28").Value = "  -2.51%  "

# Row 29
$ws.Range("E29").Value = "  +0.01%  "

# Row 30
$ws.Range("D30").Value = "3.954.99"
$ws.Range("E30").Value = "  +0.70%  "

# Row 31
$ws.Range("E31").Value = "  +7.03%  "

# Row 32
$ws.Range("E32").Value = "  +2.42%  "

# Row 33
$ws.Range("E33").Value = "  +2.82%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.58"
$ws.Range("E34").Value = "  +2.60%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.178"
$ws.Range("E35").Value = "  +7.55%  "

# Row 36
$ws.Range("E36").Value = "  +0.01%  "

# Row 37
$ws.Range("E37").Value = "  +1.23%  "

# Row 38
$ws.Range("D38").Value = "3.760.46"
$ws.Range("E38").Value = "  +0.85%  "

# Row 39
$ws.Range("E39").Value = "  +0.10%  "

# Row 40
$ws.Range("E40").Value = "  -0.66%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.94"
$ws.Range("E41").Value = "  +2.13%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.964"
$ws.Range("E42").Value = "  -0.31%  "

# Row 43
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("E44").Value = "  +10.46%  "

# Row 45
$ws.Range("E45").Value = "  +0.00%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.39"
$ws.Range("E46").Value = "  +5.49%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "158.92"
$ws.Range("E47").Value = "  +3.94%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "48.14"
$ws.Range("E48").Value = "  +3.02%  "

# Row 49
$ws.Range("E49").Value = "  +0.78%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000293"
$ws.Range("E50").Value = "  +6.11%  "

# Row 51
$ws.Range("E51").Value = "  +3.86%  "
